# Updated cryptos list with latest price/volume values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.491.34"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "1.819.40"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'316.58"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("E7").Value = "  -3.36%  "
$ws.Range("D8").Value = "'0.3881"
$ws.Range("E8").Value = "  -2.93%  "
$ws.Range("D9").Value = "'0.08471"
$ws.Range("E9").Value = "  +8.91%  "
$ws.Range("D10").Value = "'41.83"
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").Value = "'1.111"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").Value = "'6.425"
$ws.Range("D13").Value = "'21.02"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").Value = "'1.002"
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").Value = "'7.502"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").Value = "1.822.73"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").Value = "'0.00001137"
$ws.Range("E17").Value = "  +3.87%  "
$ws.Range("D18").Value = "'92.82"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").Value = "'0.06671"
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("D20").Value = "'17.71"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("D23").Value = "28.534.04"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").Value = "'11.40"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").Value = "'2.273"
$ws.Range("E25").Value = "  +1.47%  "
$ws.Range("D26").Value = "'20.99"
$ws.Range("D27").Value = "'159.33"
$ws.Range("E27").Value = "  +1.44%  "
$ws.Range("D28").Value = "2.028.79"
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("D29").Value = "'2.410"
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").Value = "'0.1086"
$ws.Range("E31").Value = "  -3.87%  "
$ws.Range("D32").Value = "'1.095"
$ws.Range("E32").Value = "  -6.11%  "
$ws.Range("D33").Value = "'5.730"
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("D34").Value = "'0.07434"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").Value = "'3.678"
$ws.Range("D36").Value = "'0.2234"
$ws.Range("E36").Value = "  -1.95%  "
$ws.Range("D37").Value = "'0.02360"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "'5.201"
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("D39").Value = "'8.764"
$ws.Range("E39").Value = "  -2.08%  "
$ws.Range("D40").Value = "'0.6322"
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("D41").Value = "'11.26"
$ws.Range("E41").Value = "  -1.58%  "
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "'3.777"
$ws.Range("E45").Value = "  +1.65%  "
$ws.Range("D46").Value = "'0.5938"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("D47").Value = "'126.20"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("D48").Value = "'1.989"
$ws.Range("E48").Value = "  -0.82%  "
$ws.Range("D49").Value = "'1.200"
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("D50").Value = "'0.06975"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").Value = "'74.33"
$ws.Range("E51").Value = "  -0.54%  "
